$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "408.63" are not
# silently converted to numbers by Excel (the source data stores these as strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '61.889.19'
$ws.Range("E2").Value = '  -1.16%  '

$ws.Range("D3").Value = '3.410.16'
$ws.Range("E3").Value = '  -0.67%  '

$ws.Range("D5").Value = '408.63'
$ws.Range("E5").Value = '  +0.32%  '

$ws.Range("D6").Value = '128.98'
$ws.Range("E6").Value = '  -1.77%  '

$ws.Range("E7").Value = '  +5.90%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  +5.30%  '

$ws.Range("E10").Value = '  +1.27%  '

$ws.Range("D11").Value = '42.68'
$ws.Range("E11").Value = '  +1.43%  '

$ws.Range("D12").Value = '0.0000216'
$ws.Range("E12").Value = '  +42.15%  '

$ws.Range("D13").Value = '9.15'
$ws.Range("E13").Value = '  +8.76%  '

$ws.Range("E14").Value = '  -0.26%  '

$ws.Range("D15").Value = '3.953.97'
$ws.Range("E15").Value = '  -0.51%  '

$ws.Range("D16").Value = '21.17'
$ws.Range("E16").Value = '  +6.69%  '

$ws.Range("D17").Value = '3.421.81'
$ws.Range("E17").Value = '  -0.57%  '

$ws.Range("D18").Value = '12.45'
$ws.Range("E18").Value = '  +7.51%  '

$ws.Range("E19").Value = '  +6.45%  '

$ws.Range("D20").Value = '61.861.89'
$ws.Range("E20").Value = '  -1.05%  '

$ws.Range("D21").Value = '448.92'
$ws.Range("E21").Value = '  +43.69%  '

$ws.Range("D22").Value = '91.23'
$ws.Range("E22").Value = '  +8.02%  '

$ws.Range("D23").Value = '3.21'
$ws.Range("E23").Value = '  +0.85%  '

$ws.Range("D24").Value = '13.11'
$ws.Range("E24").Value = '  +2.02%  '

$ws.Range("E25").Value = '  +3.36%  '

$ws.Range("D26").Value = '9.27'
$ws.Range("E26").Value = '  +14.28%  '

$ws.Range("D27").Value = '32.93'
$ws.Range("E27").Value = '  +10.77%  '

$ws.Range("E28").Value = '  +0.69%  '

$ws.Range("E29").Value = '  -2.35%  '

$ws.Range("E30").Value = '  -1.61%  '

$ws.Range("D31").Value = '12.05'
$ws.Range("E31").Value = '  +5.88%  '

$ws.Range("E32").Value = '  -1.27%  '

$ws.Range("E33").Value = '  -0.78%  '

$ws.Range("D34").Value = '42.64'
$ws.Range("E34").Value = '  -4.47%  '

$ws.Range("E35").Value = '  -0.13%  '

$ws.Range("D36").Value = '0.0500'
$ws.Range("E36").Value = '  +3.34%  '

$ws.Range("D37").Value = '53.79'
$ws.Range("E37").Value = '  +4.00%  '

$ws.Range("E39").Value = '  +2.40%  '

$ws.Range("E40").Value = '  +7.04%  '

$ws.Range("D41").Value = '2.94'
$ws.Range("E41").Value = '  -1.26%  '

$ws.Range("D42").Value = '0.318'
$ws.Range("E42").Value = '  -1.91%  '

$ws.Range("D43").Value = '142.26'
$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D44").Value = '4.24'
$ws.Range("E44").Value = '  +7.83%  '

$ws.Range("E45").Value = '  +1.03%  '

$ws.Range("E46").Value = '  +14.35%  '

$ws.Range("D47").Value = '16.57'
$ws.Range("E47").Value = '  -1.54%  '

$ws.Range("D48").Value = '22.25'
$ws.Range("E48").Value = '  +4.15%  '

$ws.Range("D49").Value = '0.143'
$ws.Range("E49").Value = '  +19.69%  '

$ws.Range("E50").Value = '  +8.59%  '

$ws.Range("D51").Value = '3.756.68'
$ws.Range("E51").Value = '  -0.66%  '

# Restore the default "Normal" style on column D so no stray number format/style
# index is left behind on cells (matches original unstyled inline-string cells).
$ws.Range("D2:D51").Style = "Normal"